$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-run of RU 1001 simulation (without crop) refreshed the bootstrap/
# simulation-derived share values below; row 7 ("My taxes ... global
# problems" (Global Nation, 2024)) no longer overflows (#NUM!) and now
# reports 0 for every country/group.

# Row 2
$ws.Range("B2").Value = 0.693991148793208
$ws.Range("L2").Value = 0.754180914635801

# Row 3
$ws.Range("B3").Value = 0.641501546421187
$ws.Range("L3").Value = 0.759131699735939

# Row 4
$ws.Range("B4").Value = 0.704960018034767
$ws.Range("D4").Value = 0.767261516731427
$ws.Range("E4").Value = 0.757541746745335
$ws.Range("F4").Value = 0.874984008964506
$ws.Range("G4").Value = 0.847206298409435
$ws.Range("H4").Value = 0.842339547896952
$ws.Range("I4").Value = 0.65799487387268
$ws.Range("J4").Value = 0.656332785129309
$ws.Range("K4").Value = 0.703874631903231
$ws.Range("L4").Value = 0.777883926828007
$ws.Range("M4").Value = 0.92816201896394
$ws.Range("N4").Value = 0.562406199574745

# Row 5
$ws.Range("D5").Value = 0.438407527923474
$ws.Range("E5").Value = 0.439022530105316
$ws.Range("F5").Value = 0.693097173376912
$ws.Range("H5").Value = 0.510821694970915
$ws.Range("I5").Value = 0.457151305629475
$ws.Range("N5").Value = 0.400676565778138

# Row 6
$ws.Range("B6").Value = 0.591975291077805
$ws.Range("D6").Value = 0.432966587857026
$ws.Range("E6").Value = 0.623667021276442
$ws.Range("F6").Value = 0.765866950290479
$ws.Range("G6").Value = 0.633729019702456
$ws.Range("H6").Value = 0.702512213956321
$ws.Range("I6").Value = 0.575699232924345
$ws.Range("J6").Value = 0.53058178548063
$ws.Range("K6").Value = 0.586217425103406
$ws.Range("L6").Value = 0.572781630922761
$ws.Range("M6").Value = 0.887006592651732
$ws.Range("N6").Value = 0.554098197443099

# Row 7
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
